$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-10, columns E-T (G..N..Q..R etc.)
# Columns: E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$data = @{
    2  = @{ E=3; G=99.883077; H=299.649231; I=0.3917580942718377; J=0.3917580942718377; K=3; M=21.09934133333334; N=63.29802400000001; O=0.2917236204149438; P=0.2917236204149438; Q=2107.467135046616; R=18967.20421541954; S=0.1142850895878393; T=0.1142850895878394 }
    3  = @{ E=3; G=99.883077; H=299.649231; I=0.3917580942718377; J=0.3917580942718377; K=3; M=35.81943766666667; N=107.458313; O=0.4952465516465762; P=0.4952465516465762; Q=3577.755650556367; R=32199.8008550073; S=0.1940168452677619; T=0.1940168452677619 }
    4  = @{ E=3; G=99.883077; H=299.649231; I=0.3917580942718377; J=0.3917580942718377; K=3; M=15.40769666666667; N=46.22309; O=0.2130298279384801; P=0.2130298279384801; Q=1538.96815254931; R=13850.71337294379; S=0.08345615941623644; T=0.08345615941623646 }
    5  = @{ E=3; G=124.0161413333333; H=372.048424; I=0.4864119993789693; J=0.4864119993789694; K=3; M=21.09934133333334; N=63.29802400000001; O=0.2917236204149438; P=0.2917236204149438; Q=2616.658896834909; R=23549.93007151418; S=0.1418978694721043; T=0.1418978694721044 }
    6  = @{ E=3; G=124.0161413333333; H=372.048424; I=0.4864119993789693; J=0.4864119993789694; K=3; M=35.81943766666667; N=107.458313; O=0.4952465516465762; P=0.4952465516465762; Q=4442.188444149857; R=39979.69599734872; S=0.2408938653719511; T=0.2408938653719512 }
    7  = @{ E=3; G=124.0161413333333; H=372.048424; I=0.4864119993789693; J=0.4864119993789694; K=3; M=15.40769666666667; N=46.22309; O=0.2130298279384801; P=0.2130298279384801; Q=1910.803087434462; R=17197.22778691016; S=0.1036202645349139; T=0.1036202645349139 }
    8  = @{ E=3; G=31.06188766666667; H=93.18566300000001; I=0.121829906349193; J=0.121829906349193; K=3; M=21.09934133333334; N=63.29802400000001; O=0.2917236204149438; P=0.2917236204149438; Q=655.3853703366569; R=5898.468333029913; S=0.03554066135500012; T=0.03554066135500013 }
    9  = @{ E=3; G=31.06188766666667; H=93.18566300000001; I=0.121829906349193; J=0.121829906349193; K=3; M=35.81943766666667; N=107.458313; O=0.4952465516465762; P=0.4952465516465762; Q=1112.619349085169; R=10013.57414176652; S=0.06033584100686314; T=0.06033584100686316 }
    10 = @{ E=3; G=31.06188766666667; H=93.18566300000001; I=0.121829906349193; J=0.121829906349193; K=3; M=15.40769666666667; N=46.22309; O=0.2130298279384801; P=0.2130298279384801; Q=478.5921430620745; R=4307.32928755867; S=0.02595340398732972; T=0.02595340398732973 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
